$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scotland Championship")

# Row 20
$ws.Range("B20").Value = 6845235
$ws.Range("E20").Value = "Inverness CT"
$ws.Range("F20").Value = "Dunfermline"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = "D"
$ws.Range("L20").Value = 2.45
$ws.Range("M20").Value = 3.1
$ws.Range("N20").Value = 2.75
$ws.Range("O20").Value = 2.375
$ws.Range("P20").Value = 3.1
$ws.Range("Q20").Value = 3.1
$ws.Range("R20").Value = -0.25
$ws.Range("S20").Value = 2.025
$ws.Range("T20").Value = 1.775
$ws.Range("U20").Value = 2.5
$ws.Range("V20").Value = 1.975
$ws.Range("W20").Value = 1.825
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = 2.1
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = -0.5
$ws.Range("AB20").Value = 0.3875
$ws.Range("AC20").Value = -1
$ws.Range("AD20").Value = 0.825

# Row 21
$ws.Range("B21").Value = 6845237
$ws.Range("E21").Value = "Raith"
$ws.Range("F21").Value = "Queens Park"
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = "H"
$ws.Range("L21").Value = 2.375
$ws.Range("M21").Value = 3.5
$ws.Range("N21").Value = 2.6
$ws.Range("O21").Value = 2.1
$ws.Range("P21").Value = 3.6
$ws.Range("Q21").Value = 3.1
$ws.Range("R21").Value = -0.25
$ws.Range("S21").Value = 1.875
$ws.Range("T21").Value = 1.925
$ws.Range("U21").Value = 2.75
$ws.Range("V21").Value = 1.975
$ws.Range("W21").Value = 1.825
$ws.Range("X21").Value = 1.1
$ws.Range("Y21").Value = -1
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = 0.875
$ws.Range("AB21").Value = -1
$ws.Range("AC21").Value = 0.9750000000000001
$ws.Range("AD21").Value = -1

# Row 22
$ws.Range("B22").Value = 6845236
$ws.Range("E22").Value = "Morton"
$ws.Range("F22").Value = "Partick"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = "A"
$ws.Range("L22").Value = 2.3
$ws.Range("M22").Value = 3.6
$ws.Range("N22").Value = 2.6
$ws.Range("O22").Value = 2.7
$ws.Range("P22").Value = 3.75
$ws.Range("Q22").Value = 2.3
$ws.Range("R22").Value = 0.25
$ws.Range("S22").Value = 1.775
$ws.Range("T22").Value = 2.025
$ws.Range("U22").Value = 2.5
$ws.Range("V22").Value = 1.8
$ws.Range("W22").Value = 2
$ws.Range("X22").Value = -1
$ws.Range("Y22").Value = -1
$ws.Range("Z22").Value = 1.3
$ws.Range("AA22").Value = -1
$ws.Range("AB22").Value = 1.025
$ws.Range("AC22").Value = 0.8
$ws.Range("AD22").Value = -1

# Row 50
$ws.Range("B50").Value = 6845274
$ws.Range("E50").Value = "Partick"
$ws.Range("F50").Value = "Ayr"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 2
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = "D"
$ws.Range("L50").Value = 1.8
$ws.Range("M50").Value = 3.5
$ws.Range("N50").Value = 4
$ws.Range("O50").Value = 1.909
$ws.Range("P50").Value = 3.6
$ws.Range("Q50").Value = 3.4
$ws.Range("R50").Value = -0.5
$ws.Range("S50").Value = 1.95
$ws.Range("T50").Value = 1.85
$ws.Range("U50").Value = 3
$ws.Range("V50").Value = 2
$ws.Range("W50").Value = 1.8
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 2.6
$ws.Range("Z50").Value = -1
$ws.Range("AA50").Value = -1
$ws.Range("AB50").Value = 0.8500000000000001
$ws.Range("AC50").Value = 1
$ws.Range("AD50").Value = -1

# Row 51
$ws.Range("B51").Value = 6845270
$ws.Range("E51").Value = "Airdrieonians"
$ws.Range("F51").Value = "Dundee Utd"
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 2
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1
$ws.Range("K51").Value = "A"
$ws.Range("L51").Value = 3.75
$ws.Range("M51").Value = 3.5
$ws.Range("N51").Value = 1.85
$ws.Range("O51").Value = 5.25
$ws.Range("P51").Value = 4
$ws.Range("Q51").Value = 1.533
$ws.Range("R51").Value = 1
$ws.Range("S51").Value = 1.9
$ws.Range("T51").Value = 1.95
$ws.Range("U51").Value = 2.75
$ws.Range("V51").Value = 1.85
$ws.Range("W51").Value = 2
$ws.Range("X51").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z51").Value = 0.5329999999999999
$ws.Range("AA51").Value = -1
$ws.Range("AB51").Value = 0.95
$ws.Range("AC51").Value = -1
$ws.Range("AD51").Value = 1

# Row 141
$ws.Range("B141").Value = 6975417
$ws.Range("E141").Value = "Airdrieonians"
$ws.Range("F141").Value = "Arbroath"
$ws.Range("G141").Value = 5
$ws.Range("H141").Value = 2
$ws.Range("I141").Value = 2
$ws.Range("J141").Value = 1
$ws.Range("K141").Value = "H"
$ws.Range("L141").Value = 1.909
$ws.Range("M141").Value = 3.4
$ws.Range("N141").Value = 3.6
$ws.Range("O141").Value = 1.45
$ws.Range("P141").Value = 4
$ws.Range("Q141").Value = 6
$ws.Range("R141").Value = -1
$ws.Range("S141").Value = 1.775
$ws.Range("T141").Value = 2.025
$ws.Range("U141").Value = 2.5
$ws.Range("V141").Value = 1.825
$ws.Range("W141").Value = 1.975
$ws.Range("X141").Value = 0.45
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = -1
$ws.Range("AA141").Value = 0.7749999999999999
$ws.Range("AB141").Value = -1
$ws.Range("AC141").Value = 0.825
$ws.Range("AD141").Value = -1

# Row 142
$ws.Range("B142").Value = 6975416
$ws.Range("E142").Value = "Inverness CT"
$ws.Range("F142").Value = "Ayr"
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = 1
$ws.Range("K142").Value = "A"
$ws.Range("L142").Value = 2.1
$ws.Range("M142").Value = 3.4
$ws.Range("N142").Value = 3.1
$ws.Range("O142").Value = 2.15
$ws.Range("P142").Value = 3.1
$ws.Range("Q142").Value = 3.2
$ws.Range("R142").Value = -0.25
$ws.Range("S142").Value = 1.9
$ws.Range("T142").Value = 1.9
$ws.Range("U142").Value = 2.25
$ws.Range("V142").Value = 1.9
$ws.Range("W142").Value = 1.9
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 2.2
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8999999999999999
$ws.Range("AC142").Value = 0.8999999999999999
$ws.Range("AD142").Value = -1

# Row 143
$ws.Range("B143").Value = 6957817
$ws.Range("E143").Value = "Partick"
$ws.Range("F143").Value = "Morton"
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = 0
$ws.Range("K143").Value = "H"
$ws.Range("L143").Value = 1.8
$ws.Range("M143").Value = 3.6
$ws.Range("N143").Value = 3.8
$ws.Range("O143").Value = 2.3
$ws.Range("P143").Value = 3.1
$ws.Range("Q143").Value = 3
$ws.Range("R143").Value = -0.25
$ws.Range("S143").Value = 2.025
$ws.Range("T143").Value = 1.775
$ws.Range("U143").Value = 2.25
$ws.Range("V143").Value = 1.9
$ws.Range("W143").Value = 1.9
$ws.Range("X143").Value = 1.3
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 1.025
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 0.8999999999999999
$ws.Range("AD143").Value = -1

# Row 144
$ws.Range("B144").Value = 6975418
$ws.Range("E144").Value = "Queens Park"
$ws.Range("F144").Value = "Raith"
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = 0
$ws.Range("K144").Value = "D"
$ws.Range("L144").Value = 3.25
$ws.Range("M144").Value = 3.25
$ws.Range("N144").Value = 2.1
$ws.Range("O144").Value = 3.2
$ws.Range("P144").Value = 3.1
$ws.Range("Q144").Value = 2.2
$ws.Range("R144").Value = 0.25
$ws.Range("S144").Value = 1.9
$ws.Range("T144").Value = 1.95
$ws.Range("U144").Value = 2.5
$ws.Range("V144").Value = 2
$ws.Range("W144").Value = 1.85
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 2.1
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.45
$ws.Range("AB144").Value = -0.5
$ws.Range("AC144").Value = -1
$ws.Range("AD144").Value = 0.8500000000000001

# Row 145
$ws.Range("B145").Value = 6975419
$ws.Range("E145").Value = "Morton"
$ws.Range("F145").Value = "Dunfermline"
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = 0
$ws.Range("K145").Value = "A"
$ws.Range("L145").Value = 1.8
$ws.Range("M145").Value = 3.5
$ws.Range("N145").Value = 3.5
$ws.Range("O145").Value = 2.1
$ws.Range("P145").Value = 3.25
$ws.Range("Q145").Value = 3
$ws.Range("R145").Value = -0.25
$ws.Range("S145").Value = 1.95
$ws.Range("T145").Value = 1.85
$ws.Range("U145").Value = 2.25
$ws.Range("V145").Value = 1.95
$ws.Range("W145").Value = 1.85
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 2
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.8500000000000001
$ws.Range("AC145").Value = -1
$ws.Range("AD145").Value = 0.8500000000000001

# Row 148
$ws.Range("B148").Value = 6994674
$ws.Range("E148").Value = "Ayr"
$ws.Range("F148").Value = "Queens Park"
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 2
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = 1
$ws.Range("K148").Value = "A"
$ws.Range("L148").Value = 2.1
$ws.Range("M148").Value = 3.4
$ws.Range("N148").Value = 2.875
$ws.Range("O148").Value = 2.25
$ws.Range("P148").Value = 3.5
$ws.Range("Q148").Value = 2.6
$ws.Range("R148").Value = -0.25
$ws.Range("S148").Value = 2.025
$ws.Range("T148").Value = 1.775
$ws.Range("U148").Value = 2.75
$ws.Range("V148").Value = 1.85
$ws.Range("W148").Value = 1.95
$ws.Range("X148").Value = -1
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = 1.6
$ws.Range("AA148").Value = -1
$ws.Range("AB148").Value = 0.7749999999999999
$ws.Range("AC148").Value = 0.425
$ws.Range("AD148").Value = -0.5

